# Adds "green hydrogen if" and "low carbon hydrogen if" subscript rows
# (rows 12 and 13) to both the Calcs sheet and the SYIEUEFbIPaF sheet,
# mirroring the existing "hydrogen if" row (row 11), and adds a blank
# formatted spacer row (row 14) below the new data on SYIEUEFbIPaF.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Calcs": copy formatting of row 11 down into rows 12 & 13,
# then fill in the labels and the B$2*0.8-style formulas.
# ---------------------------------------------------------------
$calcs = $wb.Worksheets.Item("Calcs")

$calcs.Range("A11:I11").Copy()
$calcs.Range("A12:I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$calcs.Range("A12").Value = "green hydrogen if"
$calcs.Range("A13").Value = "low carbon hydrogen if"

$cols = @("B","C","D","E","F","G","H","I")
foreach ($row in 12,13) {
    foreach ($col in $cols) {
        $calcs.Range("$col$row").Formula = "=$col`$2*0.8"
    }
}

[void]$calcs.Range("A14").Select()

# ---------------------------------------------------------------
# Sheet "SYIEUEFbIPaF": copy formatting of row 11 down into rows
# 12 & 13, fill in labels and 1/Calcs!xx formulas, then add a blank
# but formatted row 14 below.
# ---------------------------------------------------------------
$main = $wb.Worksheets.Item("SYIEUEFbIPaF")

$main.Range("A11:I11").Copy()
$main.Range("A12:I13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$main.Range("A12").Value = "green hydrogen if"
$main.Range("A13").Value = "low carbon hydrogen if"

foreach ($row in 12,13) {
    foreach ($col in $cols) {
        $main.Range("$col$row").Formula = "=1/Calcs!$col$row"
    }
}

foreach ($col in $cols) {
    $main.Range("$col" + "14").NumberFormat = $main.Range("$col" + "2").NumberFormat
}

[void]$main.Range("A14").Select()
[void]$main.Activate()
